# Fix song detail. Revert codes
#
# This script reproduces (as closely as the COM surface allows) the
# changes described by the target diff:
#   - Remove the huge placeholder row (old row 30, containing a long
#     run of backtick characters) which also removes its now-unused
#     shared string.
#   - Add a new bug entry as row 28:
#       A28 = "BUG", B28 = "Search bar focus white background",
#       C28 = "General", E28 = "TrungDQ"
#     which also adds a new shared string for the B28 text.
#   - Update the sheet selection to A29 (the cell right below the new
#     last row of data).
#
# Note: because the old placeholder string is removed from the middle
# of the shared string table, Excel automatically renumbers/reuses the
# shared string indices for every other cell that referenced strings
# after it (rows 12, 19 and 20) - this happens transparently when the
# workbook is saved, we don't need to touch those cells at all.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row holding the giant backtick placeholder text (old row 30).
$ws.Rows.Item(30).Delete()

# Add the new bug entry in row 28 (BUG / Search bar focus white background / General / - / TrungDQ).
$ws.Range("A28").Value = "BUG"
$ws.Range("B28").Value = "Search bar focus white background"
$ws.Range("C28").Value = "General"
$ws.Range("E28").Value = "TrungDQ"

# Match the updated selection shown in the workbook (cell below the new data).
$ws.Range("A29").Select()
